# khl_referees_stats_1369.xlsx — refresh of the "Главные" (main) and
# "Линейные" (linesmen) referee-stats sheets: updated per-referee season
# totals for a handful of officials, plus the as_of_utc refresh stamp
# bumped on every data row of both sheets.

$wb = $excel.ActiveWorkbook

$newStamp = "2025-11-24 03:04:16"

# Column letters -> 1-based column numbers used by Cells.Item(row, col).
# A=1 B=2 C=3 D=4 E=5 F=6 G=7 H=8 I=9 J=10 K=11 L=12 ... W=23 ... AA=27

function Update-RefereeRow {
    param($ws, $row, $values)

    foreach ($col in $values.Keys) {
        $ws.Cells.Item($row, $col).Value = $values[$col]
    }
}

# ---------------------------------------------------------------------
# Sheet "Главные"
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Главные")

Update-RefereeRow $wsMain 9  @{ 3=28; 4=434; 5=227; 6=207; 7=15.5; 8=8.109999999999999; 9=7.39; 10=111; 11=101; 23=26 }
Update-RefereeRow $wsMain 11 @{ 3=20; 4=492; 5=228; 7=24.6; 8=11.4; 9=13.2; 10=99 }
Update-RefereeRow $wsMain 15 @{ 3=18; 4=311; 5=148; 6=163; 7=17.28; 8=8.220000000000001; 9=9.06; 10=54; 11=74; 12=4; 23=10 }
Update-RefereeRow $wsMain 17 @{ 3=18; 4=284; 5=103; 6=181; 7=15.78; 8=5.72; 9=10.06; 10=49; 11=73; 23=8 }
Update-RefereeRow $wsMain 18 @{ 3=27; 4=412; 5=193; 6=219; 7=15.26; 8=7.15; 9=8.109999999999999; 10=84; 11=102; 12=3; 23=10 }
Update-RefereeRow $wsMain 21 @{ 3=24; 4=328; 5=146; 6=182; 7=13.67; 8=6.08; 9=7.58; 10=63; 11=76; 23=6 }
Update-RefereeRow $wsMain 25 @{ 3=28; 4=451; 5=218; 7=16.11; 8=7.79; 9=8.32; 10=104 }

# Refresh timestamp for every data row (2-26) of "Главные"
for ($r = 2; $r -le 26; $r++) {
    $wsMain.Cells.Item($r, 27).Value = $newStamp
}

# ---------------------------------------------------------------------
# Sheet "Линейные"
# ---------------------------------------------------------------------
$wsLine = $wb.Worksheets.Item("Линейные")

Update-RefereeRow $wsLine 11 @{ 3=17; 4=234; 5=106; 6=128; 7=13.76; 8=6.24; 9=7.53; 10=53; 11=59; 23=10 }
Update-RefereeRow $wsLine 17 @{ 3=15; 4=278; 5=155; 6=123; 7=18.53; 8=10.33; 9=8.199999999999999; 10=60; 11=54; 23=8 }
Update-RefereeRow $wsLine 19 @{ 3=25; 4=423; 5=202; 6=221; 7=16.92; 8=8.08; 9=8.84; 10=96; 11=98; 23=10 }
Update-RefereeRow $wsLine 20 @{ 3=19; 4=290; 5=149; 7=15.26; 8=7.84; 9=7.42; 10=72 }
Update-RefereeRow $wsLine 22 @{ 3=21; 4=401; 5=201; 6=200; 7=19.1; 8=9.57; 9=9.52; 10=83; 11=90; 23=26 }

# Refresh timestamp for every data row (2-26) of "Линейные"
for ($r = 2; $r -le 26; $r++) {
    $wsLine.Cells.Item($r, 27).Value = $newStamp
}
